# Auto-generated Excel COM-interop script
# Applies numeric corrections to the "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled profit-recalculation runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: ALC
$ws.Range("H18").Value = 1090
$ws.Range("I18").Value = 1090
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1090
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -806
$ws.Range("N18").ClearContents()

# Row 40: ALC
$ws.Range("H40").Value = 2312.4375
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 107: ALC
$ws.Range("H107").Value = 252
$ws.Range("I107").Value = 115
$ws.Range("K107").Value = 115
$ws.Range("M107").Value = 1805

# Row 137: ALC
$ws.Range("H137").Value = 1102.65
$ws.Range("I137").Value = 803.6
$ws.Range("J137").Value = 1999.8
$ws.Range("K137").Value = 2410.8
$ws.Range("L137").Value = 5999.4
$ws.Range("M137").Value = 139.1999999999998
$ws.Range("N137").Value = -11099.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32: ARM
$ws.Range("H32").Value = 5896.9424
$ws.Range("I32").Value = 5221.6875
$ws.Range("K32").Value = 5221.6875
$ws.Range("M32").Value = -4934.6875

# Row 88: ARM
$ws.Range("H88").Value = 3546
$ws.Range("I88").Value = 3141.2
$ws.Range("K88").Value = 3141.2
$ws.Range("M88").Value = -2735.2

# Row 91: ARM
$ws.Range("H91").Value = 3546
$ws.Range("I91").Value = 3141.2
$ws.Range("K91").Value = 3141.2
$ws.Range("M91").Value = -1737.2

# Row 110: ARM
$ws.Range("H110").Value = 5000
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -2955
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86: BSM
$ws.Range("H86").Value = 2225.56
$ws.Range("I86").Value = 2111.3333
$ws.Range("J86").Value = 2519.2856
$ws.Range("K86").Value = 2111.3333
$ws.Range("L86").Value = 2519.2856
$ws.Range("M86").Value = -988.3332999999998
$ws.Range("N86").Value = -4765.2856

# Row 89: BSM
$ws.Range("H89").Value = 2225.56
$ws.Range("I89").Value = 2111.3333
$ws.Range("J89").Value = 2519.2856
$ws.Range("K89").Value = 10556.6665
$ws.Range("L89").Value = 12596.428
$ws.Range("M89").Value = -4940.666499999999
$ws.Range("N89").Value = -23828.428

# Row 107: BSM
$ws.Range("H107").Value = 20628.572
$ws.Range("I107").Value = 2020
$ws.Range("J107").Value = 67150
$ws.Range("K107").Value = 2020
$ws.Range("L107").Value = 67150
$ws.Range("M107").Value = -100
$ws.Range("N107").Value = -70990

$ws = $wb.Worksheets.Item("CRP")
# Row 31: CRP
$ws.Range("H31").Value = 2234.9736
$ws.Range("I31").Value = 2245.5588
$ws.Range("J31").Value = 2145
$ws.Range("K31").Value = 2245.5588
$ws.Range("L31").Value = 2145
$ws.Range("M31").Value = -1950.5588
$ws.Range("N31").Value = -2735

# Row 34: CRP
$ws.Range("H34").Value = 2234.9736
$ws.Range("I34").Value = 2245.5588
$ws.Range("J34").Value = 2145
$ws.Range("K34").Value = 2245.5588
$ws.Range("L34").Value = 2145
$ws.Range("M34").Value = -2043.5588
$ws.Range("N34").Value = -2549

# Row 41: CRP
$ws.Range("H41").Value = 10009.333

# Row 50: CRP
$ws.Range("H50").Value = 7647.875
$ws.Range("J50").Value = 7647.875
$ws.Range("L50").Value = 7647.875
$ws.Range("N50").Value = -8897.875

# Row 51: CRP
$ws.Range("H51").Value = 9698.166999999999
$ws.Range("J51").Value = 8819.799999999999
$ws.Range("L51").Value = 8819.799999999999
$ws.Range("N51").Value = -10291.8

# Row 59: CRP
$ws.Range("H59").Value = 15877
$ws.Range("J59").Value = 15877
$ws.Range("L59").Value = 15877
$ws.Range("N59").Value = -18167

# Row 60: CRP
$ws.Range("H60").Value = 6053.3335
$ws.Range("I60").Value = 2796.5
$ws.Range("J60").Value = 6983.857
$ws.Range("K60").Value = 2796.5
$ws.Range("L60").Value = 6983.857
$ws.Range("M60").Value = -2285.5
$ws.Range("N60").Value = -8005.857

# Row 61: CRP
$ws.Range("H61").Value = 9698.166999999999
$ws.Range("J61").Value = 8819.799999999999
$ws.Range("L61").Value = 8819.799999999999
$ws.Range("N61").Value = -9515.799999999999

# Row 68: CRP
$ws.Range("H68").Value = 24111.111
$ws.Range("J68").Value = 24111.111
$ws.Range("L68").Value = 24111.111
$ws.Range("N68").Value = -25609.111

# Row 71: CRP
$ws.Range("H71").Value = 24111.111
$ws.Range("J71").Value = 24111.111
$ws.Range("L71").Value = 72333.333
$ws.Range("N71").Value = -79821.333

# Row 74: CRP
$ws.Range("H74").Value = 13203.5
$ws.Range("J74").Value = 13203.5
$ws.Range("L74").Value = 13203.5
$ws.Range("N74").Value = -14951.5

# Row 77: CRP
$ws.Range("H77").Value = 13203.5
$ws.Range("J77").Value = 13203.5
$ws.Range("L77").Value = 39610.5
$ws.Range("N77").Value = -48346.5

# Row 122: CRP
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("CUL")
# Row 18: CUL
$ws.Range("H18").Value = 661.875
$ws.Range("I18").Value = 333.33334
$ws.Range("J18").Value = 1084.2858
$ws.Range("K18").Value = 1000.00002
$ws.Range("L18").Value = 3252.8574
$ws.Range("M18").Value = -831.0000200000001
$ws.Range("N18").Value = -3590.8574

# Row 119: CUL
$ws.Range("H119").Value = 9500.333000000001
$ws.Range("I119").Value = 6001
$ws.Range("J119").Value = 11250
$ws.Range("K119").Value = 18003
$ws.Range("L119").Value = 33750
$ws.Range("M119").Value = -13165
$ws.Range("N119").Value = -43426

# Row 131: CUL
$ws.Range("H131").Value = 812.22
$ws.Range("J131").Value = 820.0204
$ws.Range("L131").Value = 2460.0612
$ws.Range("N131").Value = -12540.0612

$ws = $wb.Worksheets.Item("LTW")
# Row 16: LTW
$ws.Range("H16").Value = 7451.375
$ws.Range("I16").Value = 8955.538
$ws.Range("J16").Value = 933.3333
$ws.Range("K16").Value = 8955.538
$ws.Range("L16").Value = 933.3333
$ws.Range("M16").Value = -8785.538
$ws.Range("N16").Value = -1273.3333

# Row 22: LTW
$ws.Range("H22").Value = 512.25
$ws.Range("I22").Value = 529.6
$ws.Range("J22").Value = 483.33334
$ws.Range("K22").Value = 529.6
$ws.Range("L22").Value = 483.33334
$ws.Range("M22").Value = -234.6
$ws.Range("N22").Value = -1073.33334

# Row 27: LTW
$ws.Range("H27").Value = 512.25
$ws.Range("I27").Value = 529.6
$ws.Range("J27").Value = 483.33334
$ws.Range("K27").Value = 529.6
$ws.Range("L27").Value = 483.33334
$ws.Range("M27").Value = -422.6
$ws.Range("N27").Value = -697.33334

$ws = $wb.Worksheets.Item("WVR")
# Row 107: WVR
$ws.Range("H107").Value = 1080.875
$ws.Range("I107").Value = 1171.9231
$ws.Range("K107").Value = 3515.7693
$ws.Range("M107").Value = -1595.7693
